# Remove the unused "transaction_status_blockchain" column (L) from the
# payment list test fixture. This shifts the former column M (which only
# had a single stray value, "Test", in row 3) left into column L, and
# drops the now-orphaned shared strings ("transaction_status_blockchain"
# and "test_text_link111") from the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("L").Delete()
